$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers must be forced to Text
# so Excel stores them as strings (matching the source data which is inline text),
# not auto-converted to numeric values.
$textCells = @("D5","D7","D8","D9","D10","D11","D12","D13","D15","D17","D18","D19","D20","D22","D24","D25","D27","D28","D29","D30","D31","D33","D34","D35","D36","D37","D38","D39","D40","D41","D43","D44","D45","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Column D (Price) updates
$ws.Range("D2").Value = "30.144.34"
$ws.Range("D3").Value = "1.922.61"
$ws.Range("D5").Value = "331.32"
$ws.Range("D7").Value = "0.5220"
$ws.Range("D8").Value = "0.4099"
$ws.Range("D9").Value = "0.08541"
$ws.Range("D10").Value = "43.45"
$ws.Range("D11").Value = "1.130"
$ws.Range("D12").Value = "22.49"
$ws.Range("D13").Value = "6.427"
$ws.Range("D14").Value = "1.922.55"
$ws.Range("D15").Value = "7.443"
$ws.Range("D17").Value = "96.15"
$ws.Range("D18").Value = "0.00001116"
$ws.Range("D19").Value = "0.06717"
$ws.Range("D20").Value = "18.35"
$ws.Range("D22").Value = "6.048"
$ws.Range("D23").Value = "30.155.25"
$ws.Range("D24").Value = "11.35"
$ws.Range("D25").Value = "2.221"
$ws.Range("D26").Value = "2.144.42"
$ws.Range("D27").Value = "21.19"
$ws.Range("D28").Value = "159.66"
$ws.Range("D29").Value = "2.464"
$ws.Range("D30").Value = "129.38"
$ws.Range("D31").Value = "1.084"
$ws.Range("D33").Value = "6.118"
$ws.Range("D34").Value = "3.644"
$ws.Range("D35").Value = "0.02513"
$ws.Range("D36").Value = "0.06611"
$ws.Range("D37").Value = "0.2222"
$ws.Range("D38").Value = "5.244"
$ws.Range("D39").Value = "1.239"
$ws.Range("D40").Value = "8.967"
$ws.Range("D41").Value = "0.6543"
$ws.Range("D43").Value = "1.247"
$ws.Range("D44").Value = "0.6185"
$ws.Range("D45").Value = "13.26"
$ws.Range("D47").Value = "2.096"
$ws.Range("D48").Value = "1.252"
$ws.Range("D49").Value = "125.16"
$ws.Range("D50").Value = "1.164"
$ws.Range("D51").Value = "79.95"

# Column E (Volume 1h) updates
$ws.Range("E3").Value = "  +2.55%  "
$ws.Range("E4").Value = "  -0.68%  "
$ws.Range("E5").Value = "  +4.87%  "
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("E7").Value = "  +2.36%  "
$ws.Range("E8").Value = "  +5.00%  "
$ws.Range("E9").Value = "  +2.13%  "
$ws.Range("E10").Value = "  +3.86%  "
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("E12").Value = "  +9.87%  "
$ws.Range("E13").Value = "  +3.21%  "
$ws.Range("E14").Value = "  +2.35%  "
$ws.Range("E15").Value = "  +2.14%  "
$ws.Range("E17").Value = "  +5.32%  "
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("E20").Value = "  +3.42%  "
$ws.Range("E21").Value = "  -0.67%  "
$ws.Range("E22").Value = "  +2.20%  "
$ws.Range("E23").Value = "  +5.48%  "
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("E26").Value = "  +2.54%  "
$ws.Range("E27").Value = "  +2.56%  "
$ws.Range("E28").Value = "  -1.31%  "
$ws.Range("E29").Value = "  +1.73%  "
$ws.Range("E30").Value = "  +1.60%  "
$ws.Range("E31").Value = "  +3.91%  "
$ws.Range("E32").Value = "  +1.37%  "
$ws.Range("E33").Value = "  +6.44%  "
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("E35").Value = "  +2.48%  "
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("E37").Value = "  +2.59%  "
$ws.Range("E38").Value = "  +4.21%  "
$ws.Range("E39").Value = "  +4.75%  "
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("E41").Value = "  +2.43%  "
$ws.Range("E42").Value = "  +5.96%  "
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("E44").Value = "  +2.88%  "
$ws.Range("E45").Value = "  +1.36%  "
$ws.Range("E46").Value = "  +2.37%  "
$ws.Range("E47").Value = "  +4.48%  "
$ws.Range("E48").Value = "  +2.58%  "
$ws.Range("E49").Value = "  +2.46%  "
$ws.Range("E50").Value = "  +1.74%  "
$ws.Range("E51").Value = "  +4.53%  "
